$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the last used data row based on column A (data starts at row 2).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

# Update column C ("Förändrad") for every data row, changing the old date
# serial 45177 (2023-09-08) to the new serial 45178 (2023-09-09).
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 3).Value = 45178
}
